# Equipment_Database_python.xlsx — "many changes - main and selection working"
#
# 1. Update the selection on the "hammer" sheet (was G3 -> now I6).
# 2. Rename the "drilling rig" sheet to "drill rig".
# 3. Switch to / update the selection on the "drill rig" sheet (was K1 -> now F15)
#    and leave it as the active tab, matching the saved workbook's active sheet.

$wb = $excel.ActiveWorkbook

$hammer = $wb.Worksheets.Item("hammer")
$hammer.Activate()
[void]$hammer.Range("I6").Select()

$rig = $wb.Worksheets.Item("drilling rig")
$rig.Name = "drill rig"
$rig.Activate()
[void]$rig.Range("F15").Select()
